$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated values for rows 2-14 (columns B, C, D), and new row 15 (A, B, C, D)
$data = @(
    @(2010, 5358.10551236704, 1884.314774870098, 3461.334214751512),
    @(2011, 5006.665047726751, 1852.659666844038, 3172.015352190998),
    @(2012, 4651.369846296123, 1765.768588620275, 2912.798323131551),
    @(2013, 5325.417811275645, 1947.18275893511, 3388.168450604801),
    @(2014, 5216.459315523864, 1931.720870240265, 3301.474255074935),
    @(2015, 4798.380818556626, 1840.354768007607, 2989.66815042156),
    @(2016, 3957.083268959925, 1649.987421557196, 2366.296606397183),
    @(2017, 4905.534709438249, 1783.861287562288, 3105.990277509399),
    @(2018, 3674.09187657283, 1426.83803107733, 2257.82713039374),
    @(2019, 4677.649771924687, 1659.031684836401, 3007.178579354852),
    @(2020, 5054.757029258187, 1751.676621647945, 3283.467850698226),
    @(2021, 4816.546311839128, 1766.00625707088, 3050.637410347314),
    @(2022, 5081.950793909221, 1830.245176528996, 3249.596748692005),
    @(2023, 5453.00133138303, 1954.835190631959, 3498.166140751071)
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
    $row++
}
